# Case and Fatality Demographics Data Updated
# Update underlying counts on each of the six summary sheets; the
# percentage columns and "Total"/"Grand Total" rows are formulas and
# recompute automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Cases by Age Group"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B3").Value  = 1409
$ws.Range("B4").Value  = 3898
$ws.Range("B5").Value  = 15907
$ws.Range("B6").Value  = 17455
$ws.Range("B7").Value  = 15307
$ws.Range("B8").Value  = 12918
$ws.Range("B9").Value  = 4676
$ws.Range("B10").Value = 3163
$ws.Range("B11").Value = 1920
$ws.Range("B12").Value = 1271
$ws.Range("B13").Value = 1961
$ws.Activate()
$ws.Range("B2:B14").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Cases by Gender"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 27435
$ws.Range("B3").Value = 51848

# ---------------------------------------------------------------------
# Sheet 3: "Cases by RaceEthnicity"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 987
$ws.Range("B3").Value = 13142
$ws.Range("B4").Value = 28717
$ws.Range("B5").Value = 593
$ws.Range("B6").Value = 28009
$ws.Range("B7").Value = 8734
$ws.Activate()
$ws.Range("B14").Select()

# ---------------------------------------------------------------------
# Sheet 4: "Fatalities by Age Group"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B4").Value  = 35
$ws.Range("B5").Value  = 271
$ws.Range("B6").Value  = 895
$ws.Range("B7").Value  = 2603
$ws.Range("B8").Value  = 5874
$ws.Range("B9").Value  = 4849
$ws.Range("B10").Value = 6236
$ws.Range("B11").Value = 6863
$ws.Range("B12").Value = 6753
$ws.Range("B13").Value = 16914
$ws.Activate()
$ws.Range("B2:B14").Select()

# ---------------------------------------------------------------------
# Sheet 5: "Fatalities by Gender"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 21524
$ws.Range("B3").Value = 29788
$ws.Activate()
$ws.Range("B2:B4").Select()

# ---------------------------------------------------------------------
# Sheet 6: "Fatalities by Race-Ethnicity"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1092
$ws.Range("B3").Value = 5244
$ws.Range("B4").Value = 23822
$ws.Range("B5").Value = 281
$ws.Range("B6").Value = 20851

# This is the sheet that ends up active/selected once the report refresh
# is finished.
$ws.Activate()
$ws.Range("E15").Select()
